$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data, plus a couple of row content swaps (rows 39/40).
# D-column "Price" values are forced to remain as text (ClearFormats after assignment)
# because many of them are numeric-looking strings (e.g. "8.50", "0.160", "1.00") that
# Excel COM would otherwise silently coerce into numbers, losing formatting/precision.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "87.491.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.176.02"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -6.18%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -7.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "611.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.382"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -8.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.669"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.177.42"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.536"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -16.00%  "
$ws.Range("E12").Value = "  +4.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -15.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.762.44"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.28"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.344.55"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.16"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -14.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.166.18"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.03"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.46"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -10.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "416.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -10.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.50"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -12.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -11.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.19"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -7.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.90"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.338.58"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "73.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -9.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000130"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -9.32%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.160"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -17.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "544.83"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -8.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.25"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -13.06%  "
$ws.Range("E34").Value = "  -17.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.74"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -8.15%  "
$ws.Range("E36").Value = "  -13.14%  "
$ws.Range("E37").Value = "  -9.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "21.85"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -9.01%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "21.81"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("E41").Value = "  -6.78%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -12.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.370"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -15.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.02"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.07%  "
$ws.Range("E46").Value = "  -8.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.38%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  -14.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.97"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -13.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.701"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -12.21%  "
